$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1 ("Input"): replace the 17-column (A:Q) layout with the new
# 16-column (A:P) standard-template layout, re-mapping / re-deriving every
# value from the old table.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Input")

# Wipe all existing content + formatting so we can rebuild cleanly.
$ws1.Cells.Clear()

$headers1 = @("발주일자","납기일자","거래처명","거래처 이메일","납품처명","납품처 이메일","프로젝트명","대분류","중분류","소분류","품목명","규격","수량","단가","총금액","비고")
for ($c = 0; $c -lt $headers1.Length; $c++) {
    $ws1.Cells.Item(1, $c + 1).Value = $headers1[$c]
}

# Row data: 발주일자, 납기일자, 거래처명, 거래처 이메일, 납품처명, 납품처 이메일, 프로젝트명,
#           대분류, 중분류, 소분류, 품목명, 규격, 수량, 단가, 총금액, 비고(blank)
$rows1 = @(
    @("2025-08-30","2025-10-09","제이비엔지니어링","제이비엔지니어링@example.com","힐스테이트 도곡동1차","delivery@example.com","힐스테이트 도곡동1차","1. 원자재비","4) ALUM. 창호","A. 압출","5월 청구분","KS규격-1",46,4910,248446),
    @("2025-09-15","2025-09-03","제이비엔지니어링","제이비엔지니어링@example.com","힐스테이트 도곡동1차","delivery@example.com","힐스테이트 도곡동1차","5. 운반비","일반자재","기타","운반비","KS규격-2",1,0,0),
    @("2025-08-21","2025-09-06","제이비엔지니어링","제이비엔지니어링@example.com","힐스테이트 도곡동1차","delivery@example.com","힐스테이트 도곡동1차","1. 원자재비","4) ALUM. 창호","A. 압출","IJ-15861","KS규격-3",1,458040,503844),
    @("2025-09-05","2025-09-25","제이비엔지니어링","제이비엔지니어링@example.com","힐스테이트 도곡동1차","delivery@example.com","힐스테이트 도곡동1차","5. 운반비","일반자재","기타","5월 운반비","KS규격-4",1,0,0)
)

$numericCols = @(13, 14, 15)  # M=수량, N=단가, O=총금액 (1-based column index)

for ($i = 0; $i -lt $rows1.Length; $i++) {
    $r = $i + 2
    $row = $rows1[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $col = $c + 1
        $val = $row[$c]
        if ($numericCols -contains $col) {
            $ws1.Cells.Item($r, $col).Value = $val
        } else {
            # Force text storage (dates like "2025-08-30" would otherwise be
            # auto-converted to date serials) using a leading quote prefix,
            # then strip the resulting quote-prefix formatting residue.
            $ws1.Cells.Item($r, $col).Value = "'" + $val
            $ws1.Cells.Item($r, $col).ClearFormats()
        }
    }
    # Column P (비고) is left empty for every data row, matching the source.
}

# ---------------------------------------------------------------------------
# Sheets 2 & 3 ("갑지" / "을지"): drop the trailing empty 비고 (col I) cells
# on every data row so they become truly absent rather than empty strings.
# ---------------------------------------------------------------------------
foreach ($name in @("갑지", "을지")) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("I2:I5").ClearContents()
}
